# The commit deletes the "under_airtemp" and "under_rh" columns (K:L) from
# the "traits_and_envi" worksheet. Deleting these two entire columns shifts
# all subsequent columns (M:W) left by two positions (into K:U), removes the
# now-unused shared strings for those two headers, and updates the
# dimension/row spans accordingly - all handled automatically by Excel when
# entire columns are deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire columns K and L ("under_airtemp" and "under_rh"), shifting
# the remaining columns to the left.
$ws.Range("K1:L1").EntireColumn.Delete() | Out-Null

# The author's saved file shows the active selection parked at R17 (which is
# where the old "under_airtemp"/"under_rh" data used to sit, now occupied by
# the shifted-left "LMA" column) rather than the original K1.
$ws.Range("R17").Select() | Out-Null
